$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.808.16'
$ws.Range("E2").Value = '  -2.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.526.13'
$ws.Range("E3").Value = '  -1.66%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '612.46'
$ws.Range("E5").Value = '  +4.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.96'
$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.631'
$ws.Range("E7").Value = '  +1.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.214'
$ws.Range("E9").Value = '  -0.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.657'
$ws.Range("E10").Value = '  +0.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.46'
$ws.Range("E11").Value = '  -1.95%  '

$ws.Range("E12").Value = '  -3.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.68'
$ws.Range("E13").Value = '  +1.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.097.27'
$ws.Range("E14").Value = '  +2.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '615.49'
$ws.Range("E15").Value = '  +8.61%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '12.80'
$ws.Range("E16").Value = '  +3.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.925.71'
$ws.Range("E17").Value = '  -1.87%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.08'
$ws.Range("E18").Value = '  -0.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.560.73'
$ws.Range("E19").Value = '  -1.50%  '

$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.995'
$ws.Range("E21").Value = '  -1.77%  '

$ws.Range("E22").Value = '  -0.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '104.93'
$ws.Range("E23").Value = '  +10.59%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.68'
$ws.Range("E24").Value = '  +1.78%  '

$ws.Range("E25").Value = '  -1.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.02'
$ws.Range("E26").Value = '  +2.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.92'
$ws.Range("E27").Value = '  -3.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.04'
$ws.Range("E28").Value = '  +9.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.05'
$ws.Range("E29").Value = '  +4.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.03'
$ws.Range("E30").Value = '  -3.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.42'
$ws.Range("E31").Value = '  +0.80%  '

$ws.Range("E32").Value = '  +0.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '63.74'
$ws.Range("E33").Value = '  -0.64%  '

$ws.Range("E34").Value = '  +12.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '535.61'
$ws.Range("E35").Value = '  -2.28%  '

$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("E37").Value = '  -6.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.398'
$ws.Range("E38").Value = '  -5.06%  '

$ws.Range("E39").Value = '  +3.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.78'
$ws.Range("E40").Value = '  -2.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.550.97'
$ws.Range("E41").Value = '  +0.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0₃0775'
$ws.Range("E42").Value = '  -4.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.141'
$ws.Range("E43").Value = '  +3.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0461'
$ws.Range("E44").Value = '  +3.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.95'
$ws.Range("E45").Value = '  +0.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.143'
$ws.Range("E46").Value = '  +4.51%  '

$ws.Range("E47").Value = '  -3.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.95'
$ws.Range("E48").Value = '  -3.84%  '

$ws.Range("E49").Value = '  +0.49%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.26'
$ws.Range("E50").Value = '  -1.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.35'
$ws.Range("E51").Value = '  -7.55%  '
